# "Fabricated Cable IDs" (Table3) row corrections.
#
# The earth wire row (EARTH-B / EARTH) is moved up from row 13 to row 10
# (ahead of the two "HE" hotend rows and the "DC24-B1" row, which each
# shift down by one), and its AWG is corrected from 20 to 16. The other
# three rows simply carry their prior contents down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: becomes the old row 13 (EARTH-B / EARTH), AWG corrected 20 -> 16
$ws.Range("A10").Value2 = "EARTH-B"
$ws.Range("B10").Value2 = "EARTH"
$ws.Range("C10").Value2 = 16
$ws.Range("D10").Value2 = 1
$ws.Range("E10").Value2 = 60
$ws.Range("G10").Value2 = "Ring"
$ws.Range("H10").Value2 = "Ring"
$ws.Range("I10").Value2 = "PSU Earth to Frame"

# --- Row 11: becomes the old row 10 (HE-A / HE)
$ws.Range("A11").Value2 = "HE-A"
$ws.Range("B11").Value2 = "HE"
$ws.Range("C11").Value2 = 20
$ws.Range("D11").Value2 = 2
$ws.Range("E11").Value2 = 55
$ws.Range("G11").Value2 = "(None)"
$ws.Range("H11").Value2 = "MF3 PM2"
$ws.Range("I11").ClearContents()

# --- Row 12: becomes the old row 11 (HE-B / HE)
$ws.Range("A12").Value2 = "HE-B"
$ws.Range("B12").Value2 = "HE"
$ws.Range("C12").Value2 = 20
$ws.Range("D12").Value2 = 2
$ws.Range("E12").Value2 = 150
$ws.Range("G12").Value2 = "MF3 M2"
$ws.Range("H12").Value2 = "MF3 PM2"
$ws.Range("I12").ClearContents()

# --- Row 13: becomes the old row 12 (DC24-B1 / DC24)
$ws.Range("A13").Value2 = "DC24-B1"
$ws.Range("B13").Value2 = "DC24"
$ws.Range("C13").Value2 = 20
$ws.Range("D13").Value2 = 2
$ws.Range("E13").Value2 = 55
$ws.Range("G13").Value2 = "Ring"
$ws.Range("H13").Value2 = "(None)"
$ws.Range("I13").Value2 = "PSU DC to 5V Buck"

# Column F ("Column1") is the calculated Table3 column
# (=[Pin Count]*[Len (cm)]); leaving its formula alone lets it recompute
# from the new Pin Count / Len values above.

# Matches the author's final cursor position recorded in the saved file.
[void]$ws.Activate()
[void]$ws.Range("D13").Select()
